# update and run DE density calcs
#
# Adds two new worksheets, "area_lores_basic" and "area_pop_sum_basic", that
# contain the same kind of pandas describe()/summary output as the existing
# "area_lores" / "area_pop_sum" sheets, recomputed on the "basic" density
# geounits dataset for Dresden.

$wb = $excel.ActiveWorkbook

$wsLores  = $wb.Worksheets.Item("area_lores")
$wsPopSum = $wb.Worksheets.Item("area_pop_sum")

# ---------------------------------------------------------------------------
# 1) area_lores_basic - duplicate of area_lores' sheet (layout/format/labels
#    all preserved exactly), then update the numeric results column.
# ---------------------------------------------------------------------------
$wsLores.Copy([System.Reflection.Missing]::Value, $wsPopSum)
$wsLoresBasic = $wb.Worksheets.Item($wsPopSum.Index + 1)
$wsLoresBasic.Name = "area_lores_basic"

$wsLoresBasic.Range("B2").Value = 29
$wsLoresBasic.Range("B3").Value = 11.3193469022865
$wsLoresBasic.Range("B4").Value = 12.65240078940706
$wsLoresBasic.Range("B5").Value = 3.210396959446752
$wsLoresBasic.Range("B6").Value = 4.857497380800166
$wsLoresBasic.Range("B7").Value = 5.909689885751256
$wsLoresBasic.Range("B8").Value = 10.35208147337903
$wsLoresBasic.Range("B9").Value = 58.38958368320595

# ---------------------------------------------------------------------------
# 2) area_pop_sum_basic - duplicate of area_pop_sum's sheet, then update the
#    numeric results column.
# ---------------------------------------------------------------------------
$wsPopSum.Copy([System.Reflection.Missing]::Value, $wsLoresBasic)
$wsPopSumBasic = $wb.Worksheets.Item($wsLoresBasic.Index + 1)
$wsPopSumBasic.Name = "area_pop_sum_basic"

$wsPopSumBasic.Range("B1").Value = 0
$wsPopSumBasic.Range("B2").Value = 328.2610601663084
$wsPopSumBasic.Range("B3").Value = 512829
$wsPopSumBasic.Range("B4").Value = 1562.259622692327
